$d = $word.ActiveDocument

# --- 1) Remove the stray _GoBack bookmark from the title paragraph ("Progress Report 2") ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2) Split the "Week 6..." sentence and insert "(stage) " before "milestone", ---
#        leaving the cursor's new _GoBack bookmark right after the inserted text.
$find = $d.Content
$found = $find.Find.Execute("power up milestone in week 7", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse to just after "power up " (i.e. right before "milestone")
    $insertPoint = $find.Start + 9
    $target = $d.Range($insertPoint, $insertPoint)

    # Track revisions while inserting so the new text lands in its own run,
    # matching how Word keeps a freshly-typed span distinct from the
    # surrounding (unchanged) text instead of silently re-merging runs.
    $wasTracking = $d.TrackRevisions
    $d.TrackRevisions = $true
    $target.InsertAfter("(stage) ")
    $d.TrackRevisions = $wasTracking
    $d.Revisions.AcceptAll()

    # Word automatically drops the _GoBack bookmark at the point of the last edit.
    $goBackPoint = $insertPoint + 8
    $goBackRange = $d.Range($goBackPoint, $goBackPoint)
    $d.Bookmarks.Add("_GoBack", $goBackRange)
}
